$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line entries ("line7", "line8") are inserted into the "name"
# sequence right after "line6" and before the "extr*" entries. Since the
# B column cells for rows 8-15 keep referencing the same relative
# position in that sequence, their displayed text shifts down by two
# slots; rewrite them explicitly with their new text so the sheet reads
# correctly (line7, line8, extr1..extr6), then add the two brand-new
# rows 16/17 (extr7, extr8) at the end.

# Row 8: now "line7" (C,D updated; E unchanged TRUE)
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

# Row 9: now "line8" (C updated; E flips to TRUE)
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

# Row 10: now "extr1" (C,D updated; E flips to TRUE)
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# Row 11: now "extr2" (C,D updated; E flips to TRUE)
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

# Row 12: now "extr3" (C updated; E flips to FALSE)
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10
$ws.Range("E12").Value = $false

# Row 13: now "extr4" (D updated; E flips to TRUE)
$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# Row 14: now "extr5" (C,D updated; E flips to TRUE)
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# Row 15: now "extr6" (C,D updated; E stays FALSE)
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# New row 16 ("extr7") - copy row 15's A-column formatting (bold + border)
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# New row 17 ("extr8")
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false

$excel.CutCopyMode = 0
